# ============================================================================
# Edit: update ZBP_02_zasazeni_domacnosti workbook with a new survey wave
#   - Sheet "data"   (AN = "25. 1. 2021" wave -> corrected to "25. 1. 2022";
#                      AO = brand-new "22. 2. 2022" wave column)
#   - Sheet "pocetR" (AM = "25. 1. 2021" wave -> corrected to "25. 1. 2022";
#                      AN = brand-new "22. 2. 2022" wave column)
# ============================================================================

$wb = $excel.ActiveWorkbook
$wsData   = $wb.Worksheets.Item("data")
$wsPocetR = $wb.Worksheets.Item("pocetR")

# ----------------------------------------------------------------------
# Sheet "data": column AN = 40, column AO = 41
# Row 1 holds the wave-date headers, rows 2-67 hold the indicator values.
# ----------------------------------------------------------------------
$sheet1Data = @(
    ,@(1, "25. 1. 2022", "22. 2. 2022")
    ,@(2, 0.16, 0.16)
    ,@(3, 0.15, 0.15)
    ,@(4, 0.6899999999999999, 0.6899999999999999)
    ,@(5, 0.09, 0.07000000000000001)
    ,@(6, 0.15, 0.16)
    ,@(7, 0.76, 0.77)
    ,@(8, 0.07000000000000001, 0.03)
    ,@(9, 0.12, 0.13)
    ,@(10, 0.8100000000000001, 0.84)
    ,@(11, 0.24, 0.25)
    ,@(12, 0.15, 0.14)
    ,@(13, 0.61, 0.61)
    ,@(14, 0.23, 0.27)
    ,@(15, 0.16, 0.12)
    ,@(16, 0.61, 0.61)
    ,@(17, 0.08, 0.06)
    ,@(18, 0.15, 0.16)
    ,@(19, 0.77, 0.78)
    ,@(20, 0.06, 0.03)
    ,@(21, 0.11, 0.12)
    ,@(22, 0.83, 0.85)
    ,@(23, 0.26, 0.28)
    ,@(24, 0.12, 0.12)
    ,@(25, 0.62, 0.6)
    ,@(26, 0.23, 0.25)
    ,@(27, 0.15, 0.14)
    ,@(28, 0.62, 0.61)
    ,@(29, 0.28, 0.27)
    ,@(30, 0.2, 0.2)
    ,@(31, 0.52, 0.53)
    ,@(32, 0.1, 0.11)
    ,@(33, 0.13, 0.12)
    ,@(34, 0.77, 0.77)
    ,@(35, 0.03, 0.03)
    ,@(36, 0.07000000000000001, 0.07000000000000001)
    ,@(37, 0.9, 0.9)
    ,@(38, 0.18, 0.18)
    ,@(39, 0.17, 0.16)
    ,@(40, 0.65, 0.66)
    ,@(41, 0.15, 0.15)
    ,@(42, 0.14, 0.14)
    ,@(43, 0.71, 0.71)
    ,@(44, 0.67, 0.68)
    ,@(45, 0.11, 0.1)
    ,@(46, 0.22, 0.22)
    ,@(47, 0.31, 0.28)
    ,@(48, 0.42, 0.41)
    ,@(49, 0.27, 0.31)
    ,@(50, 0.06, 0.06)
    ,@(51, 0.08, 0.08)
    ,@(52, 0.86, 0.86)
    ,@(53, 0.13, 0.1)
    ,@(54, 0.14, 0.17)
    ,@(55, 0.73, 0.73)
    ,@(56, 0, 0)
    ,@(57, 0.13, 0.14)
    ,@(58, 0.87, 0.86)
    ,@(59, 0.12, 0.06)
    ,@(60, 0.2, 0.25)
    ,@(61, 0.68, 0.6899999999999999)
    ,@(62, 0.02, 0.05)
    ,@(63, 0.1, 0.08)
    ,@(64, 0.88, 0.87)
    ,@(65, 0.07000000000000001, 0.04)
    ,@(66, 0.14, 0.08)
    ,@(67, 0.79, 0.88)
)

for ($i = 0; $i -lt $sheet1Data.Count; $i++) {
    $entry = $sheet1Data[$i]
    $r  = $entry[0]
    $an = $entry[1]
    $ao = $entry[2]
    $wsData.Cells.Item($r, 40).Value = $an
    $wsData.Cells.Item($r, 41).Value = $ao
}

# Footer caption in row 68 mentions the latest update date.
$wsData.Cells.Item(68, 1).Value = "Život během pandemie, Zasažení domácností, % respondentů celkově a ve skupinách, aktualizace 2. 3. 2022"

# ----------------------------------------------------------------------
# Sheet "pocetR": column AM = 39, column AN = 40
# Row 1 holds the wave-date headers, rows 2-23 hold the sample sizes.
# ----------------------------------------------------------------------
$sheet2Data = @(
    ,@(1, "25. 1. 2022", "22. 2. 2022")
    ,@(2, 1582, 1510)
    ,@(3, 751, 706)
    ,@(4, 124, 123)
    ,@(5, 492, 476)
    ,@(6, 215, 205)
    ,@(7, 716, 670)
    ,@(8, 114, 113)
    ,@(9, 103, 97)
    ,@(10, 649, 630)
    ,@(11, 734, 706)
    ,@(12, 542, 517)
    ,@(13, 306, 287)
    ,@(14, 439, 409)
    ,@(15, 1143, 1101)
    ,@(16, 144, 142)
    ,@(17, 312, 288)
    ,@(18, 1126, 1080)
    ,@(19, 272, 256)
    ,@(20, 86, 86)
    ,@(21, 257, 243)
    ,@(22, 143, 129)
    ,@(23, 82, 81)
)

for ($i = 0; $i -lt $sheet2Data.Count; $i++) {
    $entry = $sheet2Data[$i]
    $r  = $entry[0]
    $am = $entry[1]
    $an = $entry[2]
    $wsPocetR.Cells.Item($r, 39).Value = $am
    $wsPocetR.Cells.Item($r, 40).Value = $an
}

# Footer caption in row 24 mentions the latest update date.
$wsPocetR.Cells.Item(24, 1).Value = "Život během pandemie, Zasažení domácností, velikost dotázaného souboru celkově a ve skupinách, aktualizace 2. 3. 2022"

# Row 24 is otherwise a row of "touched" blank filler cells spanning every
# column up to the last data column; replicate that for the new AN24 cell
# so it materializes (still empty/no value) like its neighbors.
$wsPocetR.Cells.Item(24, 40).Borders.LineStyle = 0

Write-Output "Update applied: data!AN:AO (rows 1-68), pocetR!AM:AN (rows 1-24)"
